$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week date range) ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Crime-statistics table updates (rows 15-28, 31) ---
$ws.Range("C29").Copy($ws.Range("C15"))
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 18
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = 38.461538461538
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -81.818181818181
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -25
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 18.75
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 22.222222222222
$ws.Range("L17").Value = 37.5
$ws.Range("M17").Value = 450
$ws.Range("N17").Value = -42.105263157894
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = -52
$ws.Range("L18").Value = -14.285714285714
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = -83.333333333333
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -10.63829787234
$ws.Range("I19").Value = 45
$ws.Range("J19").Value = 55
$ws.Range("K19").Value = -18.181818181818
$ws.Range("L19").Value = -19.642857142857
$ws.Range("M19").Value = 55.172413793103
$ws.Range("N19").Value = -50.54945054945
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -85.714285714285
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = -85.714285714285
$ws.Range("M20").Value = -80
$ws.Range("N20").Value = -98.305084745762
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -10.78431372549
$ws.Range("I21").Value = 101
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = -15.833333333333
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 48.529411764705
$ws.Range("N21").Value = -72.022160664819
$ws.Range("C29").Copy($ws.Range("C22"))
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -60
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -42.857142857142
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = -60
$ws.Range("C23").Value = 2
$ws.Range("C29").Copy($ws.Range("D23"))
$ws.Range("E29").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 250
$ws.Range("I23").Value = 7
$ws.Range("K23").Value = 250
$ws.Range("L23").Value = 75
$ws.Range("M23").Value = 600
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -10.714285714285
$ws.Range("I24").Value = 112
$ws.Range("J24").Value = 121
$ws.Range("K24").Value = -7.438016528925
$ws.Range("L24").Value = -13.846153846153
$ws.Range("M24").Value = -20.567375886524
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 18
$ws.Range("F25").Value = 82
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = 3.79746835443
$ws.Range("I25").Value = 91
$ws.Range("J25").Value = 85
$ws.Range("K25").Value = 7.058823529411
$ws.Range("L25").Value = -23.529411764705
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 25
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -6.666666666666
$ws.Range("I26").Value = 29
$ws.Range("J26").Value = 33
$ws.Range("K26").Value = -12.121212121212
$ws.Range("L26").Value = -3.333333333333
$ws.Range("M26").Value = -12.121212121212
$ws.Range("C29").Copy($ws.Range("C27"))
$ws.Range("C29").Copy($ws.Range("C28"))
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -14.285714285714
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -14.285714285714
$ws.Range("C39").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("K39").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("C39").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("K39").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100
$ws.Range("C39").Copy($ws.Range("J31"))
$ws.Range("J31").Value = 1
$ws.Range("K39").Copy($ws.Range("K31"))
$ws.Range("K31").Value = -100
